$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "tes"
$ws.Range("B2").Value = "set"
$ws.Range("A3").Value = "tas"
$ws.Range("B3").Value = "sat"
